# Rewrite the "question" text in column B:
#   "Use the risk of bias tool to evaluate the risk of bias concerning <X> and provide evidence for supporting it."
# becomes
#   "Evaluate the risk of bias concerning <X>."
#
# Applies to every data row (rows 2..188) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$prefix = "Use the risk of bias tool to evaluate the risk of bias concerning "
$suffix = " and provide evidence for supporting it."

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 188 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $text = $cell.Value()
    if ($text -ne $null -and $text.StartsWith($prefix) -and $text.EndsWith($suffix)) {
        $middle = $text.Substring($prefix.Length, $text.Length - $prefix.Length - $suffix.Length)
        $cell.Value = "Evaluate the risk of bias concerning " + $middle + "."
    }
}
